$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71; this shifts all rows 71.. down by one
# (row 155 content becomes row 156, dimension grows to A1:R156 automatically)
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new price record
$ws.Cells.Item(71, 1).Value = 4
$ws.Cells.Item(71, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(71, 3).Value = "Los Lagos"
$ws.Cells.Item(71, 4).Value = 44483
$ws.Cells.Item(71, 4).Style = $ws.Cells.Item(72, 4).Style
$ws.Cells.Item(71, 4).NumberFormat = $ws.Cells.Item(72, 4).NumberFormat
$ws.Cells.Item(71, 5).Value = 10
$ws.Cells.Item(71, 6).Value = 100112003
$ws.Cells.Item(71, 7).Value = "Ajo"
$ws.Cells.Item(71, 8).Value = "Chino"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 120
$ws.Cells.Item(71, 11).Value = 17500
$ws.Cells.Item(71, 12).Value = 19000
$ws.Cells.Item(71, 13).Value = 18250
$ws.Cells.Item(71, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(71, 15).Value = "China"
$ws.Cells.Item(71, 16).Value = 1825
$ws.Cells.Item(71, 17).Value = 10
$ws.Cells.Item(71, 18).Value = "Hortaliza"
